$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two patient records (rows 2 and 3) were reordered: sort the data
# block by Patient ID (column A) descending, so the row that used to be
# row 3 (Patient ID 100000226615) becomes row 2, and vice versa.
$ws.Range("A1:AC3").Sort($ws.Range("A2:A3"), 2, $null, $null, 1, $null, 1, 1)

# The mailto: hyperlink lives on the "E-Mail" column (X) of the row that
# now sits at row 2, so move the hyperlink definition from X3 to X2.
$ws.Range("X3").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("X2"), "mailto:Leaalexander@yahoo.com")
# Restore the cell's existing "Hyperlink" style (Add() stamps a fresh,
# visually-identical style otherwise).
$ws.Range("X2").Style = "Hipertaut"

# Re-select the whole of row 2 (matches the saved selection state).
$ws.Range("A2:XFD2").Select()
